# SectorGroup.xlsx update:
# The codeforiati:group-code column (previously last, column G) is moved to be
# the first of the four "category/group" columns (now column D), shifting the
# previous category-name / category-code / group-name columns one place to
# the right (D->E, E->F, F->G).
#
# This affects both the header row (handled automatically by Excel since the
# header text simply follows the column it sits in) and every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

# Header row (row 1): rotate the 4 header labels the same way as the data.
$headerRange = $ws.Range("D1:G1")
$headerRange.NumberFormat = "@"
$header = $headerRange.Value()
$groupCode = $header[1,4]
$catName   = $header[1,1]
$catCode   = $header[1,2]
$groupName = $header[1,3]
$ws.Range("D1").Value = $groupCode
$ws.Range("E1").Value = $catName
$ws.Range("F1").Value = $catCode
$ws.Range("G1").Value = $groupName

# Data rows (row 2 .. lastRow): rotate D,E,F,G -> new D = old G, new E = old D,
# new F = old E, new G = old F.
if ($lastRow -ge 2) {
    $dataRange = $ws.Range("D2:G" + $lastRow)

    # Force text formatting so numeric-looking codes (e.g. "110") keep being
    # stored as text/shared-strings rather than turning into numbers, matching
    # the original (and target) cell typing.
    $dataRange.NumberFormat = "@"

    $data = $dataRange.Value()

    $rowCount = $lastRow - 2 + 1
    $newData = New-Object 'object[,]' $rowCount,4

    for ($i = 1; $i -le $rowCount; $i++) {
        $oldCatName   = $data[$i,1]
        $oldCatCode   = $data[$i,2]
        $oldGroupName = $data[$i,3]
        $oldGroupCode = $data[$i,4]

        $newData[$i-1,0] = $oldGroupCode
        $newData[$i-1,1] = $oldCatName
        $newData[$i-1,2] = $oldCatCode
        $newData[$i-1,3] = $oldGroupName
    }

    $dataRange.Value = $newData
}
